$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 corresponds to the "Exp" class. Implementing IComparable on it
# flips the "IComparable" column's three sub-columns (ThisClass,
# AbstractInterface, and the non-generic IComparable column) from
# "TODO" (Incorrecto / red) to "Oui" (Bueno / green), matching the rest
# of the "Oui"-styled cells in that row.
$ws.Range("B5").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("G5").Value = "Oui"
$ws.Range("I5").Value = "Oui"
$ws.Range("K5").Value = "Oui"
$excel.CutCopyMode = 0

# Update the active cell selection left in the worksheet view.
$ws.Range("I19").Select()
